$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC05's "Files" query (cell B4) had its trailing LIMIT clause commented out.
# Uncomment it (turn "--LIMIT 100;" into "LIMIT 100;") to fix the query.
$current = $ws.Range("B4").Value2
$fixed = $current.Replace('--LIMIT 100;', 'LIMIT 100;')
$ws.Range("B4").Value2 = $fixed

# Leave the cursor on B3 (matches the saved selection state in the fixed workbook).
$ws.Range("B3").Select()
